# The workbook tracks a "Förändrad" (changed) date in column C for every
# data row. This update bumps that date by one day (2023-10-03 -> 2023-10-04,
# i.e. Excel serial 45202 -> 45203) for every row that currently has it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 472 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45203
